$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-31 hold the per-team data (col B = team abbreviation, col C = metric value).
# The team order and the metric values were regenerated; column A (index) is unchanged.
$ws.Cells.Item(2, 2).Value = "NOK"
$ws.Cells.Item(2, 3).Value = 13.22
$ws.Cells.Item(3, 2).Value = "POR"
$ws.Cells.Item(3, 3).Value = 12.35
$ws.Cells.Item(4, 2).Value = "NJN"
$ws.Cells.Item(4, 3).Value = 8.835714285714285
$ws.Cells.Item(5, 2).Value = "CLE"
$ws.Cells.Item(5, 3).Value = 11.87333333333333
$ws.Cells.Item(6, 2).Value = "DAL"
$ws.Cells.Item(6, 3).Value = 12.72
$ws.Cells.Item(7, 2).Value = "ATL"
$ws.Cells.Item(7, 3).Value = 11.78823529411765
$ws.Cells.Item(8, 2).Value = "SEA"
$ws.Cells.Item(8, 3).Value = 11.39375
$ws.Cells.Item(9, 2).Value = "CHA"
$ws.Cells.Item(9, 3).Value = 12.73333333333333
$ws.Cells.Item(10, 2).Value = "WAS"
$ws.Cells.Item(10, 3).Value = 12.25
$ws.Cells.Item(11, 2).Value = "MIL"
$ws.Cells.Item(11, 3).Value = 10.46428571428571
$ws.Cells.Item(12, 2).Value = "LAC"
$ws.Cells.Item(12, 3).Value = 8.862499999999999
$ws.Cells.Item(13, 2).Value = "SAS"
$ws.Cells.Item(13, 3).Value = 14.52142857142857
$ws.Cells.Item(14, 2).Value = "DET"
$ws.Cells.Item(14, 3).Value = 14.43571428571429
$ws.Cells.Item(15, 2).Value = "ORL"
$ws.Cells.Item(15, 3).Value = 13.26666666666667
$ws.Cells.Item(16, 2).Value = "UTA"
$ws.Cells.Item(16, 3).Value = 11.65714285714286
$ws.Cells.Item(17, 2).Value = "MEM"
$ws.Cells.Item(17, 3).Value = 14.79285714285714
$ws.Cells.Item(18, 2).Value = "HOU"
$ws.Cells.Item(18, 3).Value = 13.49230769230769
$ws.Cells.Item(19, 2).Value = "DEN"
$ws.Cells.Item(19, 3).Value = 12.00833333333333
$ws.Cells.Item(20, 2).Value = "LAL"
$ws.Cells.Item(20, 3).Value = 13.7
$ws.Cells.Item(21, 2).Value = "GSW"
$ws.Cells.Item(21, 3).Value = 12.89166666666667
$ws.Cells.Item(22, 2).Value = "IND"
$ws.Cells.Item(22, 3).Value = 11.92727272727273
$ws.Cells.Item(23, 2).Value = "CHI"
$ws.Cells.Item(23, 3).Value = 12.76428571428572
$ws.Cells.Item(24, 2).Value = "PHI"
$ws.Cells.Item(24, 3).Value = 12.01538461538462
$ws.Cells.Item(25, 2).Value = "BOS"
$ws.Cells.Item(25, 3).Value = 12.3875
$ws.Cells.Item(26, 2).Value = "TOR"
$ws.Cells.Item(26, 3).Value = 14.82142857142857
$ws.Cells.Item(27, 2).Value = "MIA"
$ws.Cells.Item(27, 3).Value = 13.16
$ws.Cells.Item(28, 2).Value = "SAC"
$ws.Cells.Item(28, 3).Value = 10.67857142857143
$ws.Cells.Item(29, 2).Value = "PHO"
$ws.Cells.Item(29, 3).Value = 13.20714285714286
$ws.Cells.Item(30, 2).Value = "NYK"
$ws.Cells.Item(30, 3).Value = 12.17333333333333
$ws.Cells.Item(31, 2).Value = "MIN"
$ws.Cells.Item(31, 3).Value = 11.62857142857143
